# tuned lambda threshold, this didnt make that much of a difference
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target column (E) for rows 54 and 55 changes from "> 140/80" to "> 160/100"
$ws.Range("E54").Value = "> 160/100"
$ws.Range("E55").Value = "> 160/100"

# Make room for two new data rows (56:57), then three more blank rows
# further down (at what becomes 60:62), matching the row layout growth
# from A1:P59 to A1:P64.
$ws.Rows("56:57").Insert()
$ws.Rows("60:62").Insert()

# Restore the correct row heights for the freshly inserted rows.
$ws.Rows("56:57").RowHeight = 18.75
$ws.Rows("60:62").RowHeight = 19.5

# New row 56 (a run like the others above, with a bad result)
$ws.Range("A56").Value = "ukb51139_subset.csv"
$ws.Range("B56").Value = "28012 x 1081"
$ws.Range("C56").Value = "all"
$ws.Range("D56").Value = "no events"
$ws.Range("E56").Value = "> 160/100"
$ws.Range("F56").Value = "zscore"
$ws.Range("G56").Value = "median"
$ws.Range("H56").Value = "none"
$ws.Range("I56").Value = 250
$ws.Range("L56").Value = "bad"
$ws.Range("M56").Value = "bad"
$ws.Range("N56").Value = 227
$ws.Range("N56").NumberFormat = "@"
$ws.Range("O56").NumberFormat = "@"
$ws.Range("O56").Value = " 47.5"
$ws.Range("P56").Value = "change lambda > 0.001"

# New row 57 (another run with the tuned lambda threshold)
$ws.Range("A57").Value = "ukb51139_subset.csv"
$ws.Range("B57").Value = "28012 x 1081"
$ws.Range("C57").Value = "all"
$ws.Range("D57").Value = "no events"
$ws.Range("E57").Value = "> 160/100"
$ws.Range("F57").Value = "zscore"
$ws.Range("G57").Value = "median"
$ws.Range("H57").Value = "none"
$ws.Range("I57").Value = 250
$ws.Range("L57").Value = "102.8 & 101.6"
$ws.Range("M57").Value = "80 & 80.9"
$ws.Range("N57").Value = 17
$ws.Range("O57").Value = 5.09
$ws.Range("P57").Value = "change lambda > 0.01"

# Inserting rows copies formatting down from the row above, which leaves
# I/K/N/O/P on rows 56:57 one style group "off" (and N56 was forced to a
# text number format above to preserve the " 47.5" leading space as text).
# Re-stamp the correct number formats from the row-54/55 template cells.
$ws.Range("I54:P54").Copy() | Out-Null
$ws.Range("I56").PasteSpecial(-4122) | Out-Null
$ws.Range("I55:P55").Copy() | Out-Null
$ws.Range("I57").PasteSpecial(-4122) | Out-Null
